$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
"no_zero", "oth", "oth", "oth", "oth", "oth", "oth", "no_zero", "mineral", "no_zero", "oth", "oth", "mineral", "no_zero", "oth", "oth", "oth", "no_zero", "oth", "oth", "mineraloth", "no_zero", "oth", "oth", "no_zero", "no_zero", "oth", "mineral", "oth", "mineral", "no_zero", "oth", "oth", "no_zero", "mineraloth", "oth", "oth", "oth", "no_zero", "oth", "oth", "oth", "no_zero", "oth", "no_zero", "oth", "no_zero", "oth", "mineral", "mineraloth", "oth", "no_zero", "no_zero", "no_zero", "mineraloth", "no_zero", "mineral", "mineral", "mineraloth", "mineral", "oth", "oth", "mineral", "oth", "oth", "oth", "oth", "oth", "oth", "mineral", "oth", "oth", "mineral", "oth", "oth", "oth", "no_zero", "no_zero", "no_zero", "oth", "no_zero", "oth", "no_zero", "oth", "oth", "oth", "oth", "oth", "oth", "oth", "no_zero", "oth", "oth", "mineraloth", "oth", "oth", "no_zero", "oth", "mineral", "mineral", "mineral", "oth", "mineral", "oth", "oth", "oth", "mineral", "mineraloth", "oth", "oth", "oth", "oth", "mineraloth", "oth", "oth", "mineraloth", "oth", "oth", "oth", "oth", "no_zero", "no_zero", "oth", "oth", "oth", "mineraloth", "oth", "no_zero", "oth", "oth", "oth", "oth", "oth", "oth", "oth", "oth", "no_zero", "oth", "no_zero", "oth", "oth", "oth", "no_zero", "oth", "oth", "oth", "oth", "oth", "no_zero", "no_zero", "no_zero", "oth", "oth", "no_zero", "no_zero", "oth", "oth", "mineraloth", "no_zero", "no_zero", "mineral", "no_zero", "oth", "oth", "oth", "no_zero", "oth", "mineraloth", "mineral", "no_zero", "oth", "mineraloth", "mineral", "no_zero", "no_zero", "no_zero", "oth", "oth", "no_zero", "oth", "no_zero", "oth", "oth", "no_zero", "oth", "no_zero", "no_zero", "no_zero", "oth", "oth", "no_zero", "no_zero", "no_zero", "no_zero", "oth", "mineraloth", "oth", "oth", "no_zero", "no_zero", "oth", "no_zero", "no_zero", "oth", "oth", "oth", "oth", "oth", "oth", "oth", "oth", "oth", "oth", "oth", "no_zero", "mineraloth", "no_zero", "oth", "no_zero", "mineraloth", "oth", "no_zero", "mineraloth", "no_zero", "oth", "oth", "oth", "no_zero", "oth", "oth", "no_zero", "no_zero", "no_zero", "oth", "no_zero", "oth", "oth", "oth", "oth", "oth", "mineraloth", "no_zero", "oth", "no_zero", "mineral", "no_zero", "no_zero", "no_zero", "no_zero", "oth", "oth", "oth", "oth", "oth", "oth", "no_zero", "oth", "oth", "no_zero", "oth", "oth", "mineral", "oth", "oth", "oth", "oth", "oth", "oth", "no_zero", "no_zero", "oth", "oth", "oth", "oth", "mineraloth", "oth", "mineraloth", "no_zero", "oth", "oth", "no_zero", "mineraloth", "mineraloth", "oth", "no_zero", "oth", "mineraloth", "no_zero", "oth", "oth", "mineraloth", "mineraloth", "no_zero", "no_zero", "no_zero", "no_zero", "oth", "no_zero", "oth", "no_zero", "mineraloth", "oth", "mineral", "oth", "no_zero", "no_zero", "oth", "no_zero", "no_zero", "oth", "oth", "no_zero", "oth", "no_zero", "mineraloth", "no_zero", "oth", "oth", "no_zero", "oth", "mineraloth", "mineraloth", "oth", "no_zero", "oth", "no_zero", "mineraloth", "oth", "oth", "no_zero", "no_zero", "oth", "oth", "oth", "no_zero", "oth", "oth", "oth", "oth", "no_zero", "no_zero", "no_zero", "no_zero", "mineraloth", "no_zero", "no_zero", "oth", "oth", "no_zero", "no_zero", "oth", "oth", "oth", "oth", "no_zero", "oth", "mineraloth", "no_zero", "oth", "oth", "no_zero", "oth", "oth", "oth", "oth", "oth", "oth", "oth", "oth", "mineraloth", "oth", "oth", "oth", "oth", "no_zero", "oth", "oth", "mineral", "no_zero", "oth", "no_zero", "oth", "no_zero", "mineraloth", "no_zero", "mineraloth", "no_zero", "oth", "no_zero", "oth", "oth", "no_zero", "no_zero", "mineraloth", "mineraloth", "oth", "oth", "mineraloth", "mineraloth", "oth", "oth", "oth", "oth", "no_zero", "oth", "oth", "no_zero", "oth", "oth", "no_zero", "oth", "oth", "oth", "mineral", "mineraloth", "oth", "oth", "oth", "oth", "oth", "mineraloth", "oth", "oth", "no_zero", "oth", "no_zero", "oth", "mineraloth", "oth", "no_zero", "oth", "no_zero", "no_zero", "no_zero", "oth", "no_zero", "oth", "no_zero", "no_zero", "oth", "oth", "oth", "oth", "mineraloth", "mineraloth", "no_zero", "oth", "oth", "no_zero", "no_zero", "mineraloth", "mineraloth", "oth", "oth", "oth", "oth", "oth", "no_zero", "oth", "oth", "oth", "oth", "oth", "oth", "oth", "mineral"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 33).Value = $values[$i]
}

Write-Output "done: $($values.Length) rows updated"